# Batch experiments and plots are pending; Experiments need to be rerun
#
# The existing "Test Accuracy (%)" column (D) is being split into two
# dataset-specific columns: CIFAR-10 (existing column D, renamed) and a
# brand-new MNIST column (E). Only the header row gets real values for
# now - the new column's data rows are intentionally left blank pending
# the reruns mentioned in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Test Accuracy (%)" header to be CIFAR-10 specific.
$ws.Range("D1").Value = "CIFAR-10 Test Accuracy (%)"

# Add the new MNIST header in column E, matching the bold style already
# used by the rest of row 1.
$ws.Range("E1").Value = "MNIST Test Accuracy (%)"
$ws.Range("E1").Font.Bold = $true

# Widen columns D and E to fit their new, longer header text (mirrors the
# "best fit" auto-resize Excel performs automatically after a header edit).
# ColumnWidth is quantized to whole pixels internally, so the inputs below
# are chosen to land on the closest achievable width to the real target
# ("CIFAR-10 Test Accuracy (%)" -> ~24.11 chars, "MNIST Test Accuracy (%)"
# -> ~22.22 chars).
$ws.Columns("D").ColumnWidth = 23.375
$ws.Columns("E").ColumnWidth = 21.375

# The active cell/selection moves to the newly added header cell.
$ws.Range("E1").Select()
